$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player table data (Oyuncu Adı / Pozisyon / Takım)
$data = @(
    @("Josh Hart",           "SG,SF,PF", "New York Knicks"),
    @("Kelly Oubre Jr.",     "SG,SF",    "Philadelphia 76ers"),
    @("Michael Porter Jr.",  "SF,PF",    "Denver Nuggets"),
    @("Kel'el Ware",         "C",        "Miami Heat"),
    @("Malik Beasley",       "SG,SF",    "Detroit Pistons"),
    @("Santi Aldama",        "PF,C",     "Memphis Grizzlies"),
    @("Alperen Sengün",      "C",        "Houston Rockets"),
    @("Kristaps Porzingis",  "PF,C",     "Boston Celtics"),
    @("Dyson Daniels",       "PG,SG,SF", "Atlanta Hawks"),
    @("Donovan Mitchell",    "PG,SG",    "Cleveland Cavaliers"),
    @("De'Andre Hunter",     "SF,PF",    "Atlanta Hawks"),
    @("Andrew Wiggins",      "SF,PF",    "Golden State Warriors"),
    @("Domantas Sabonis",    "C",        "Sacramento Kings"),
    @("Victor Wembanyama",   "C",        "San Antonio Spurs"),
    @("Jaden McDaniels",     "SF,PF",    "Minnesota Timberwolves"),
    @("Cam Thomas",          "SG,SF",    "Brooklyn Nets"),
    @("Donte DiVincenzo",    "PG,SG,SF", "Minnesota Timberwolves")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}
